$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '69.967.43'
$ws.Range('E2').Value = '  -0.26%  '

# Row 3
$ws.Range('D3').Value = '3.541.37'
$ws.Range('E3').Value = '  +0.63%  '

# Row 4
$ws.Range('E4').Value = '  +0.00%  '

# Row 5
$ws.Range('D5').Value = "'603.65"
$ws.Range('E5').Value = '  -0.83%  '

# Row 6
$ws.Range('D6').Value = "'195.65"
$ws.Range('E6').Value = '  -1.70%  '

# Row 7
$ws.Range('E7').Value = '  -0.66%  '

# Row 8
$ws.Range('D8').Value = "'1.00"
$ws.Range('E8').Value = '  +0.05%  '

# Row 9
$ws.Range('E9').Value = '  -4.62%  '

# Row 10
$ws.Range('D10').Value = "'0.649"
$ws.Range('E10').Value = '  -1.69%  '

# Row 11
$ws.Range('D11').Value = "'53.57"
$ws.Range('E11').Value = '  -1.41%  '

# Row 12
$ws.Range('D12').Value = "'0.0000304"
$ws.Range('E12').Value = '  -0.68%  '

# Row 13
$ws.Range('E13').Value = '  -1.84%  '

# Row 14
$ws.Range('D14').Value = '4.100.01'
$ws.Range('E14').Value = '  +0.72%  '

# Row 15
$ws.Range('D15').Value = "'594.89"
$ws.Range('E15').Value = '  -0.98%  '

# Row 16
$ws.Range('D16').Value = "'12.83"
$ws.Range('E16').Value = '  +0.65%  '

# Row 17
$ws.Range('B17').Value = 'Chainlink'
$ws.Range('C17').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D17').Value = "'19.18"
$ws.Range('E17').Value = '  +0.48%  '

# Row 18
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '70.033.07'
$ws.Range('E18').Value = '  -0.16%  '

# Row 19
$ws.Range('D19').Value = '3.516.13'
$ws.Range('E19').Value = '  +0.38%  '

# Row 20
$ws.Range('D20').Value = "'0.122"
$ws.Range('E20').Value = '  +1.75%  '

# Row 21
$ws.Range('D21').Value = "'0.987"
$ws.Range('E21').Value = '  -1.12%  '

# Row 22
$ws.Range('D22').Value = "'17.79"
$ws.Range('E22').Value = '  -0.68%  '

# Row 23
$ws.Range('B23').Value = 'Toncoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D23').Value = "'5.15"
$ws.Range('E23').Value = '  +1.44%  '

# Row 24
$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').Value = "'102.75"
$ws.Range('E24').Value = '  -1.62%  '

# Row 25
$ws.Range('E25').Value = '  -1.17%  '

# Row 26
$ws.Range('E26').Value = '  -1.62%  '

# Row 27
$ws.Range('D27').Value = "'10.76"
$ws.Range('E27').Value = '  -2.17%  '

# Row 28
$ws.Range('D28').Value = "'9.54"
$ws.Range('E28').Value = '  -3.84%  '

# Row 29
$ws.Range('D29').Value = "'33.17"
$ws.Range('E29').Value = '  -2.55%  '

# Row 30
$ws.Range('D30').Value = "'4.36"
$ws.Range('E30').Value = '  -3.76%  '

# Row 31
$ws.Range('D31').Value = "'7.06"
$ws.Range('E31').Value = '  -2.15%  '

# Row 32
$ws.Range('D32').Value = "'12.34"
$ws.Range('E32').Value = '  -2.83%  '

# Row 33
$ws.Range('D33').Value = "'0.116"
$ws.Range('E33').Value = '  -0.60%  '

# Row 34
$ws.Range('D34').Value = "'63.43"
$ws.Range('E34').Value = '  -0.99%  '

# Row 35
$ws.Range('D35').Value = '3.837.35'
$ws.Range('E35').Value = '  +3.95%  '

# Row 36
$ws.Range('E36').Value = '  +6.61%  '

# Row 37
$ws.Range('E37').Value = '  +2.61%  '

# Row 38
$ws.Range('E38').Value = '  +0.20%  '

# Row 39
$ws.Range('D39').Value = "'512.49"
$ws.Range('E39').Value = '  -2.14%  '

# Row 40
$ws.Range('D40').Value = "'0.391"
$ws.Range('E40').Value = '  -0.44%  '

# Row 41
$ws.Range('D41').Value = "'3.58"
$ws.Range('E41').Value = '  -0.21%  '

# Row 42
$ws.Range('D42').Value = "'36.48"
$ws.Range('E42').Value = '  -1.62%  '

# Row 43
$ws.Range('D43').Value = "'0.134"
$ws.Range('E43').Value = '  -2.59%  '

# Row 44
$ws.Range('D44').Value = "'0.0450"
$ws.Range('E44').Value = '  -2.84%  '

# Row 45
$ws.Range('B45').Value = 'ThetaToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D45').Value = "'2.82"
$ws.Range('E45').Value = '  -1.96%  '

# Row 46
$ws.Range('B46').Value = 'Stellar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D46').Value = "'0.140"
$ws.Range('E46').Value = '  -0.57%  '

# Row 47
$ws.Range('D47').Value = "'3.29"
$ws.Range('E47').Value = '  -0.52%  '

# Row 48
$ws.Range('B48').Value = 'THORChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D48').Value = "'8.51"
$ws.Range('E48').Value = '  -3.37%  '

# Row 49
$ws.Range('B49').Value = 'FirstDigitalUSD'
$ws.Range('C49').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D49').Value = "'1.00"
$ws.Range('E49').Value = '  +0.09%  '

# Row 50
$ws.Range('D50').Value = "'0.000246"
$ws.Range('E50').Value = '  +2.54%  '

# Row 51
$ws.Range('D51').Value = "'1.31"
$ws.Range('E51').Value = '  +3.56%  '

